# Updated cryptos list (Price / Volume(1h) columns) to match the latest
# scrape. Price values that look like plain decimal numbers are forced to
# text (NumberFormat "@") before assignment so Excel doesn't coerce them
# into numeric cells (which would silently drop significant trailing
# zeros, e.g. "86.00" -> 86). Percentage strings already contain
# surrounding spaces/"%" so they remain text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.616.22'
$ws.Range("E2").Value = '  +1.68%  '
$ws.Range("D3").Value = '1.600.11'
$ws.Range("E3").Value = '  +1.37%  '
$ws.Range("E4").Value = '  +0.47%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.28'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.515'
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '26.87'
$ws.Range("E8").Value = '  +3.85%  '
$ws.Range("E9").Value = '  +1.36%  '
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("E11").Value = '  +1.11%  '
$ws.Range("D12").Value = '1.828.43'
$ws.Range("D13").Value = '1.601.07'
$ws.Range("E13").Value = '  -0.08%  '
$ws.Range("D14").Value = '29.622.46'
$ws.Range("E15").Value = '  +3.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.74'
$ws.Range("E16").Value = '  +1.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.94'
$ws.Range("E17").Value = '  +2.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '241.72'
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("E19").Value = '  +2.62%  '
$ws.Range("E20").Value = '  +0.33%  '
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.23'
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.78'
$ws.Range("E25").Value = '  +0.98%  '
$ws.Range("E26").Value = '  +1.63%  '
$ws.Range("E27").Value = '  +0.50%  '
$ws.Range("E28").Value = '  +1.21%  '
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("E31").Value = '  +0.18%  '
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.15'
$ws.Range("E33").Value = '  +3.11%  '
$ws.Range("D34").Value = '1.426.05'
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("E35").Value = '  +2.19%  '
$ws.Range("E36").Value = '  +4.13%  '
$ws.Range("E37").Value = '  -1.99%  '
$ws.Range("E38").Value = '  +0.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0167'
$ws.Range("E39").Value = '  +2.56%  '
$ws.Range("E40").Value = '  +3.04%  '
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("E42").Value = '  +5.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '54.06'
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.807'
$ws.Range("E44").Value = '  +2.32%  '
$ws.Range("E45").Value = '  +0.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '66.19'
$ws.Range("E47").Value = '  +2.75%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.30'
$ws.Range("E48").Value = '  -0.70%  '
$ws.Range("D49").Value = '1.740.50'
$ws.Range("E49").Value = '  +1.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.00'
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("E51").Value = '  +5.89%  '
